# Apply WorldServer.xlsx "Property" sheet data-fill edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: fill in the actual server record (ID, ServerID, Name, MaxOnline, CpuCount, IP, Port) ---
# Columns with text values need the "text" (@) number format, matching the style
# already used by the other text columns (A2/B2) in this row.
$ws.Range("A2").Value = "WorldServer_1"
$ws.Range("B2").Value = "000103001"

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "WorldServer_1"

$ws.Range("D2").Value = 5000
$ws.Range("E2").Value = 1

$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "127.0.0.1"

$ws.Range("G2").Value = 3001

# --- Column widths ---
# (iron_native quantizes ColumnWidth to whole "pixel" steps of 1/7 character unit,
#  so the values below are chosen to land as close as possible to the target widths.)
$ws.Columns.Item(2).ColumnWidth = 13.714285714285714   # -> width 14.375
$ws.Columns.Item(3).ColumnWidth = 17.857142857142858   # -> width 18.5
$ws.Columns.Item(5).ColumnWidth = 11.285714285714286   # -> width 12
$ws.Columns.Item(6).ColumnWidth = 15.0                 # -> width 15.75
$ws.Columns.Item(7).ColumnWidth = 10.285714285714286   # -> width 11

# --- Selection moved to G1 ---
$ws.Range("G1").Select() | Out-Null

# --- Remove the data validations that used to live on column F ---
$ws.Range("F1").Validation.Delete()
$ws.Range("F2:F1048576").Validation.Delete()
